$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume(1h) figures
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.642.24'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.79%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.585.97'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.56%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.63'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.66%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.53'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.48%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.584.37'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +1.60%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.56%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.97%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.21'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +4.41%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.391'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.76%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.196.42'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.54%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.12'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.77%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.90%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.587.06'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.47%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.03%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '65.755.01'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.77%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.01'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -3.18%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.61'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.52%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.89'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '396.53'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.08%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +3.37%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.731.83'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.59%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.75%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000118'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +3.13%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.12'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +5.82%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.65'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +31.28%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.40'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +5.32%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.61'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +5.25%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.28%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.589.66'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.23%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '24.49'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.05%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.00%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.40'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +8.74%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +4.48%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.07'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.00%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '169.15'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.26%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0835'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +4.78%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.842'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.32%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.88'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +3.47%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.27'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +8.19%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '43.13'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.87%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.54'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.99%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.13%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.61%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.04'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +3.57%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.453.46'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +2.61%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '317.14'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +5.12%  '
